$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (R1:T1) - shared strings added in order: Red Color, Blue Color, Green Color
$ws.Range("R1").Value = "Red Color"
$ws.Range("S1").Value = "Blue Color"
$ws.Range("T1").Value = "Green Color"

# Row 2 - Venus: average RGB color formulas
$ws.Range("R2").Formula = "=(139+187+221+239)/4/255"
$ws.Range("S2").Formula = "=(145+183+216+239)/4/255"
$ws.Range("T2").Formula = "=(161+171+212+239)/4/255"

# Row 3 - Earth: Altitude Input changed from 500 to 1000, plus new color formulas
$ws.Range("Q3").Value = 1000
$ws.Range("R3").Formula = "=(5+11+227+161+178+127)/6/255"
$ws.Range("S3").Formula = "=(51+158+197+110+189+139)/6/255"
$ws.Range("T3").Formula = "=(85+210+117+71+91+59)/6/255"

# Row 4 - Jupiter: new color formulas
$ws.Range("R4").Formula = "=(64+167+210+211+144+200)/6/255"
$ws.Range("S4").Formula = "=(68+156+207+156+97+139)/6/255"
$ws.Range("T4").Formula = "=(54+134+218+126+77+58)/6/255"

# Row 5 - Uranus: new color formulas
$ws.Range("R5").Formula = "=(213+187+147+101)/4/255"
$ws.Range("S5").Formula = "=(251+225+184+134)/4/255"
$ws.Range("T5").Formula = "=(252+228+190+139)/4/255"

# Update the active selection to S1, matching the saved workbook state
$ws.Range("S1").Select() | Out-Null
